$wb = $excel.ActiveWorkbook

# --- Sheet: BFS ---
$ws = $wb.Worksheets.Item("BFS")
$ws.Range("B2").Value = 23
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = "[23, 19, 16, 12, 8]"
$ws.Range("E2").Value = 453
$ws.Range("F2").Value = 16
$ws.Range("G2").Value = 1.3125
$ws.Range("H2").Value = 0.0003414154052734375
$ws.Range("B3").Value = 41
$ws.Range("C3").Value = 18
$ws.Range("D3").Value = "[41, 38, 31, 20, 19, 18]"
$ws.Range("E3").Value = 596
$ws.Range("F3").Value = 21
$ws.Range("G3").Value = 1.19047619047619
$ws.Range("H3").Value = 0.0001759529113769531
$ws.Range("B4").Value = 32
$ws.Range("C4").Value = 39
$ws.Range("D4").Value = "[32, 31, 38, 39]"
$ws.Range("E4").Value = 253
$ws.Range("F4").Value = 15
$ws.Range("G4").Value = 1.466666666666667
$ws.Range("H4").Value = 0.0001504421234130859
$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "[12, 8, 5, 4]"
$ws.Range("E5").Value = 448
$ws.Range("F5").Value = 15
$ws.Range("G5").Value = 1.533333333333333
$ws.Range("H5").Value = 0.0001313686370849609
$ws.Range("B6").Value = 21
$ws.Range("C6").Value = 34
$ws.Range("D6").Value = "[21, 18, 14, 10, 7, 8, 9, 28, 34]"
$ws.Range("E6").Value = 1036
$ws.Range("F6").Value = 39
$ws.Range("G6").Value = 1.051282051282051
$ws.Range("H6").Value = 0.0001745223999023438

# --- Sheet: DFS ---
$ws = $wb.Worksheets.Item("DFS")
$ws.Range("B2").Value = 23
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = "[23, 19, 16, 12, 8]"
$ws.Range("E2").Value = 453
$ws.Range("F2").Value = 3835261
$ws.Range("G2").Value = 0.9999997392615522
$ws.Range("H2").Value = 1.71847128868103
$ws.Range("B3").Value = 41
$ws.Range("C3").Value = 18
$ws.Range("D3").Value = "[41, 32, 24, 23, 22, 21, 18]"
$ws.Range("E3").Value = 571
$ws.Range("F3").Value = 5113534
$ws.Range("G3").Value = 0.9999998044405298
$ws.Range("H3").Value = 2.365155696868896
$ws.Range("B4").Value = 32
$ws.Range("C4").Value = 39
$ws.Range("D4").Value = "[32, 41, 38, 39]"
$ws.Range("E4").Value = 252
$ws.Range("F4").Value = 5110434
$ws.Range("G4").Value = 0.9999998043219029
$ws.Range("H4").Value = 2.33404803276062
$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "[12, 8, 7, 4]"
$ws.Range("E5").Value = 438
$ws.Range("F5").Value = 3573789
$ws.Range("G5").Value = 0.9999997201849353
$ws.Range("H5").Value = 1.922887086868286
$ws.Range("B6").Value = 21
$ws.Range("C6").Value = 34
$ws.Range("D6").Value = "[21, 18, 14, 15, 16, 17, 30, 29, 28, 34]"
$ws.Range("E6").Value = 706
$ws.Range("F6").Value = 16513894
$ws.Range("G6").Value = 0.9999999394449305
$ws.Range("H6").Value = 4.870802879333496

# --- Sheet: BCU ---
$ws = $wb.Worksheets.Item("BCU")
$ws.Range("B2").Value = 23
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = "[23, 19, 16, 12, 8]"
$ws.Range("E2").Value = 453
$ws.Range("F2").Value = 26
$ws.Range("G2").Value = 0.7352941176470589
$ws.Range("H2").Value = 0.0002031326293945312
$ws.Range("B3").Value = 41
$ws.Range("C3").Value = 18
$ws.Range("D3").Value = "[41, 32, 24, 23, 22, 21, 18]"
$ws.Range("E3").Value = 571
$ws.Range("F3").Value = 33
$ws.Range("G3").Value = 0.8205128205128205
$ws.Range("H3").Value = 0.0002779960632324219
$ws.Range("B4").Value = 32
$ws.Range("C4").Value = 39
$ws.Range("D4").Value = "[32, 41, 38, 39]"
$ws.Range("E4").Value = 252
$ws.Range("F4").Value = 13
$ws.Range("G4").Value = 0.6
$ws.Range("H4").Value = 0.0001344680786132812
$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "[12, 8, 7, 4]"
$ws.Range("E5").Value = 438
$ws.Range("F5").Value = 29
$ws.Range("G5").Value = 0.7
$ws.Range("H5").Value = 0.0001776218414306641
$ws.Range("B6").Value = 21
$ws.Range("C6").Value = 34
$ws.Range("D6").Value = "[21, 18, 14, 15, 16, 17, 30, 29, 28, 34]"
$ws.Range("E6").Value = 706
$ws.Range("F6").Value = 34
$ws.Range("G6").Value = 0.8461538461538461
$ws.Range("H6").Value = 0.0001685619354248047

# --- Sheet: A_Estrela_Euclidiano ---
$ws = $wb.Worksheets.Item("A_Estrela_Euclidiano")
$ws.Range("B2").Value = 23
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = "[23, 19, 16, 12, 8]"
$ws.Range("E2").Value = 453
$ws.Range("F2").Value = 21
$ws.Range("G2").Value = 3.071428571428572
$ws.Range("H2").Value = 0.0002729892730712891
$ws.Range("B3").Value = 41
$ws.Range("C3").Value = 18
$ws.Range("D3").Value = "[41, 32, 24, 23, 22, 21, 18]"
$ws.Range("E3").Value = 571
$ws.Range("F3").Value = 26
$ws.Range("G3").Value = 3.071428571428572
$ws.Range("H3").Value = 0.0002608299255371094
$ws.Range("B4").Value = 32
$ws.Range("C4").Value = 39
$ws.Range("D4").Value = "[32, 41, 38, 39]"
$ws.Range("E4").Value = 252
$ws.Range("F4").Value = 11
$ws.Range("G4").Value = 3.071428571428572
$ws.Range("H4").Value = 0.0001816749572753906
$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "[12, 8, 7, 4]"
$ws.Range("E5").Value = 438
$ws.Range("F5").Value = 26
$ws.Range("G5").Value = 3.071428571428572
$ws.Range("H5").Value = 0.0002675056457519531
$ws.Range("B6").Value = 21
$ws.Range("C6").Value = 34
$ws.Range("D6").Value = "[21, 18, 14, 15, 16, 17, 30, 29, 28, 34]"
$ws.Range("E6").Value = 706
$ws.Range("F6").Value = 31
$ws.Range("G6").Value = 3.071428571428572
$ws.Range("H6").Value = 0.0002770423889160156

# --- Sheet: A_Estrela_Haversiano ---
$ws = $wb.Worksheets.Item("A_Estrela_Haversiano")
$ws.Range("B2").Value = 23
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = "[23, 19, 16, 12, 8]"
$ws.Range("E2").Value = 453
$ws.Range("F2").Value = 21
$ws.Range("G2").Value = 3.071428571428572
$ws.Range("H2").Value = 0.0003378391265869141
$ws.Range("B3").Value = 41
$ws.Range("C3").Value = 18
$ws.Range("D3").Value = "[41, 32, 24, 23, 22, 21, 18]"
$ws.Range("E3").Value = 571
$ws.Range("F3").Value = 26
$ws.Range("G3").Value = 3.071428571428572
$ws.Range("H3").Value = 0.0003247261047363281
$ws.Range("B4").Value = 32
$ws.Range("C4").Value = 39
$ws.Range("D4").Value = "[32, 41, 38, 39]"
$ws.Range("E4").Value = 252
$ws.Range("F4").Value = 11
$ws.Range("G4").Value = 3.071428571428572
$ws.Range("H4").Value = 0.0002036094665527344
$ws.Range("B5").Value = 12
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "[12, 8, 7, 4]"
$ws.Range("E5").Value = 438
$ws.Range("F5").Value = 26
$ws.Range("G5").Value = 3.071428571428572
$ws.Range("H5").Value = 0.00031280517578125
$ws.Range("B6").Value = 21
$ws.Range("C6").Value = 34
$ws.Range("D6").Value = "[21, 18, 14, 15, 16, 17, 30, 29, 28, 34]"
$ws.Range("E6").Value = 706
$ws.Range("F6").Value = 31
$ws.Range("G6").Value = 3.071428571428572
$ws.Range("H6").Value = 0.0003390312194824219
